$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = new price (optional); E = new volume text }
$updates = @{
    2  = @{ D = "59.455.30";   E = "  +0.48%  " }
    3  = @{ D = "2.600.30";    E = "  +0.44%  " }
    4  = @{               E = "  +0.09%  " }
    5  = @{ D = "538.04";      E = "  +2.91%  " }
    6  = @{ D = "141.30";      E = "  +1.23%  " }
    7  = @{               E = "  +0.14%  " }
    9  = @{               E = "  -0.70%  " }
    10 = @{               E = "  +1.47%  " }
    11 = @{ D = "0.334";       E = "  +1.61%  " }
    12 = @{               E = "  -0.95%  " }
    13 = @{ D = "3.061.66";    E = "  +0.50%  " }
    14 = @{ D = "59.360.43";   E = "  +0.72%  " }
    15 = @{ D = "20.80";       E = "  +1.22%  " }
    16 = @{ D = "2.621.15";    E = "  +1.14%  " }
    17 = @{               E = "  +0.30%  " }
    18 = @{ D = "340.86";      E = "  +0.74%  " }
    19 = @{ D = "4.36";        E = "  +1.29%  " }
    20 = @{ D = "10.09";       E = "  +0.02%  " }
    21 = @{ D = "6.35";        E = "  -1.94%  " }
    22 = @{               E = "  +0.00%  " }
    23 = @{ D = "67.36";       E = "  +1.71%  " }
    24 = @{               E = "  +1.18%  " }
    25 = @{               E = "  -1.81%  " }
    26 = @{               E = "  +0.05%  " }
    27 = @{               E = "  +2.58%  " }
    28 = @{ D = "0.0₃0742";    E = "  +2.44%  " }
    31 = @{ D = "5.81";        E = "  -1.40%  " }
    32 = @{ D = "18.79";       E = "  +0.54%  " }
    33 = @{ D = "149.85";      E = "  +0.50%  " }
    34 = @{               E = "  -0.23%  " }
    35 = @{               E = "  -0.52%  " }
    36 = @{ D = "0.838";       E = "  +2.32%  " }
    37 = @{               E = "  -0.63%  " }
    38 = @{ D = "0.823";       E = "  -0.21%  " }
    39 = @{               E = "  +0.39%  " }
    40 = @{               E = "  +0.09%  " }
    41 = @{ D = "271.53";      E = "  -0.25%  " }
    42 = @{ D = "0.599";       E = "  +1.48%  " }
    43 = @{ D = "10.75";       E = "  +0.11%  " }
    44 = @{               E = "  -0.01%  " }
    45 = @{ D = "0.0524";      E = "  +1.44%  " }
    46 = @{ D = "18.59";       E = "  +3.61%  " }
    47 = @{               E = "  +1.63%  " }
    48 = @{ D = "1.939.76";    E = "  -1.19%  " }
    49 = @{ D = "4.48";        E = "  -0.76%  " }
    50 = @{ D = "112.14";      E = "  -0.94%  " }
    51 = @{ D = "4.78";        E = "  +1.17%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals["D"]
    }
    $ws.Range("E$row").Value = $vals["E"]
}
